# anillo 2 y 3
# Aportes de elian. anillos en cuatitlan y tultitlan con enlaces para conectarse a satelite

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Anillo 2): new site info contributed by elian - tec edomex / ford cuatitlan
$ws.Range("B3").Value = "tec edomex"
$ws.Range("C3").Value = "ford cuatitlan"

# H2 previously held the placeholder text "lo llame anillo x" (shared string reused).
# It is now replaced by just "x", and new "y" / "z" values are added for rows 3 and 4.
$ws.Range("H2").Value = "x"
$ws.Range("H3").Value = "y"
$ws.Range("H4").Value = "z"

# Row 4 (Anillo 3): new site info - kio tultitlan
$ws.Range("B4").Value = "kio tultitlan"

# Update the active selection as recorded in the saved view
$ws.Range("C6").Select()
